# Updates the cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed
# quotes, matching a GitHub Actions scheduled data refresh.
#
# The Price column stores values as literal text (e.g. "26.755.82", a
# thousands-dotted price with no real decimal meaning) rather than numbers, so a
# plain Range.Value assignment of a numeric-looking string would get silently
# auto-coerced to a true number by Excel (losing both the text type and any
# trailing/format-significant digits). To keep these cells text, we briefly force
# the cell's NumberFormat to Text ("@") before writing the value, then restore the
# cell's original Style afterward so no visible formatting changes are left behind.
# The Volume(1h) column already contains a "%" sign so it's never auto-converted
# and can be written directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.755.82"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  -0.01%  "

$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.533.28"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -2.05%  "

$ws.Range("E4").Value = "  -0.04%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.22"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("E6").Value = "  -1.07%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -1.14%  "

$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.22"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -3.37%  "

$ws.Range("E10").Value = "  -0.77%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0854"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -0.80%  "

$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.751.41"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -2.00%  "

$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.530.53"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  -2.18%  "

$ws.Range("E14").Value = "  -1.54%  "

$ws.Range("E15").Value = "  -1.65%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.747.90"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -0.22%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.80"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -1.09%  "

$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "212.50"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("E19").Value = "  -1.98%  "

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0679"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +0.50%  "

$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("E22").Value = "  -2.37%  "

$ws.Range("E23").Value = "  -2.77%  "

$ws.Range("E24").Value = "  -3.57%  "

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.36"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -1.12%  "

$ws.Range("E26").Value = "  -2.75%  "

$ws.Range("E27").Value = "  -0.94%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("E29").Value = "  -1.36%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.09"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -1.24%  "

$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0453"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -2.20%  "

$ws.Range("E32").Value = "  +2.28%  "

$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.360.18"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -1.82%  "

$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("E35").Value = "  -2.64%  "

$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.954"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +3.33%  "

$ws.Range("E37").Value = "  -0.40%  "

$ws.Range("E38").Value = "  +0.41%  "

$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.518"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -0.84%  "

$ws.Range("E40").Value = "  -1.90%  "

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.72"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +6.80%  "

$ws.Range("E42").Value = "  +0.30%  "

$ws.Range("E43").Value = "  +0.44%  "

$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.45"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -1.22%  "

$ws.Range("E45").Value = "  -2.80%  "

$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.665.79"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -1.98%  "

$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.09"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -0.45%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₇0970"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -1.58%  "

$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0940"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -1.06%  "
